# Append: 2025-11-21 12:35 JST
# A new scrape run: one new listing is prepended (after the header row and the
# still-open first listing), the remaining previously-seen listings shift down
# by one row, and every still-present row's "fetched at" timestamp (col A) is
# refreshed to the new run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-11-21 12:35:38"

# --- 1. Drop the existing hyperlink relationships up front -----------------
# Row insertion shifts cell VALUES down a row but this host does not also
# renumber the worksheet's <hyperlinks> ref map, so stale links would end up
# pointing at the wrong row. Clear them now and re-add correct ones once all
# the row content is final.
$ws.Hyperlinks.Delete()

# --- 2. Insert a new row right under row 2, pushing rows 3-7 to 4-8 --------
$ws.Rows.Item(3).Insert()

# --- 3. Refresh the "fetched at" timestamp for every surviving row ---------
$ws.Range("A2").Value = $newTimestamp
$ws.Range("A4").Value = $newTimestamp
$ws.Range("A5").Value = $newTimestamp
$ws.Range("A6").Value = $newTimestamp
$ws.Range("A7").Value = $newTimestamp
$ws.Range("A8").Value = $newTimestamp

# --- 4. Populate the freshly inserted row 3 with the new listing ----------
$ws.Range("A3").Value = $newTimestamp
$ws.Range("B3").Value = "初回 Pythonのテキストエディターに機能追加依頼"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "1,000 ~ 5,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5439127"
$ws.Range("G3").Value = 190
$ws.Range("H3").Value = "🔥Python"

# --- 5. Re-create every hyperlink against the now-correct rows -------------
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5438171")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5439127")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5438740")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5438554")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5438369")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5438567")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5438092")

# --- 6. Column D narrows from 32 to 28 raw OOXML chars (COM ColumnWidth is
# offset from the stored <col width> by ~0.83, so 27.17 round-trips to 28) --
$ws.Columns.Item(4).ColumnWidth = 27.17
